# Loan RBI, Variable Instalments
# On the "Repayment schedule" sheet, a new column is inserted before the
# existing "Late" column (old column N) to make room for an additional
# instalment-schedule column. Everything from the old column N onward
# (Late / heading / Outstanding) shifts one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N; existing N/O/P data shifts to O/P/Q.
$ws.Columns("N").Insert()

# The new column picks up the width of the column immediately to its
# left (M), same as Excel's default "insert" formatting behaviour.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active/selected sheet and leave the
# cursor on K16, matching the saved view state.
$ws.Activate()
$ws.Range("K16").Select() | Out-Null
